$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(5).Insert()
$ws.Range("A993").NumberFormat = "General"
Write-Output ("UsedRange: " + $ws.UsedRange.Address())
